$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 28; $row++) {
    # Column H - PERIOD TO EXPIRE: decreases by 1 for every row
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value2 = $hCell.Value2 - 1

    # Column I - LAST UPDATE: move the progress date forward a day.
    # Prefix with an apostrophe so Excel stores it as literal text
    # (matching the existing text cells) instead of re-parsing the
    # "dd-mmm-yyyy" look-alike string into a date serial number.
    $iCell = $ws.Cells.Item($row, 9)
    $iCell.Value2 = "'04-Nov-2025"
}
